# Form the consolidated report: update the "Absent" (column H) values
# for the rows whose attendance status changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
